$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-converted to a number by Excel;
# mark them as Text before assigning so they stay strings (matching the original
# inlineStr cells), then clear the temporary formatting so the style stays General.

$ws.Range("D2").Value = '36.212.51'
$ws.Range("E2").Value = '  +2.37%  '

$ws.Range("D3").Value = '2.018.64'
$ws.Range("E3").Value = '  +6.45%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.76'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.657'
$ws.Range("E6").Value = '  -4.78%  '
$ws.Range("D6").ClearFormats()

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.93'
$ws.Range("E8").Value = '  +4.56%  '
$ws.Range("D8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '58.27'
$ws.Range("E9").Value = '  +4.23%  '
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.362'
$ws.Range("E10").Value = '  +1.40%  '
$ws.Range("D10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0714'
$ws.Range("E11").Value = '  -4.51%  '
$ws.Range("D11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0987'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.61'
$ws.Range("E13").Value = '  +4.48%  '
$ws.Range("D13").ClearFormats()

$ws.Range("D14").Value = '2.320.09'
$ws.Range("E14").Value = '  +6.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.807'
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = '2.022.68'
$ws.Range("E16").Value = '  +8.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.88'
$ws.Range("E17").Value = '  -2.65%  '
$ws.Range("D17").ClearFormats()

$ws.Range("D18").Value = '36.563.98'
$ws.Range("E18").Value = '  +3.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.85'
$ws.Range("E19").Value = '  -3.65%  '
$ws.Range("D19").ClearFormats()

$ws.Range("D20").Value = '0.0₃0818'
$ws.Range("E20").Value = '  -1.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.95'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '234.01'
$ws.Range("E22").Value = '  -4.33%  '
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.94'
$ws.Range("E23").Value = '  -5.51%  '
$ws.Range("D23").ClearFormats()

$ws.Range("E24").Value = '  -0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  -8.04%  '
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.10'
$ws.Range("E26").Value = '  -2.60%  '
$ws.Range("D26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.99'
$ws.Range("E27").Value = '  -8.64%  '
$ws.Range("D27").ClearFormats()

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.60'
$ws.Range("E28").Value = '  +6.94%  '
$ws.Range("D28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.50'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.121'
$ws.Range("E30").Value = '  -5.04%  '
$ws.Range("D30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.38'
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.73'
$ws.Range("E32").Value = '  +53.14%  '
$ws.Range("D32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0587'
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("D33").ClearFormats()

$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("E35").Value = '  -0.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.00'
$ws.Range("E36").Value = '  -5.36%  '
$ws.Range("D36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0825'
$ws.Range("E37").Value = '  +15.40%  '
$ws.Range("D37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.15'
$ws.Range("E38").Value = '  +9.55%  '
$ws.Range("D38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.838'
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.35'
$ws.Range("E40").Value = '  -7.97%  '
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '96.28'
$ws.Range("E41").Value = '  -2.63%  '
$ws.Range("D41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0213'
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("D42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.08'
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = '2.75'
$ws.Range("E44").Value = '  +13.65%  '
$ws.Range("D44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '15.89'
$ws.Range("E45").Value = '  -6.60%  '
$ws.Range("D45").ClearFormats()

$ws.Range("D46").Value = '1.315.60'
$ws.Range("E46").Value = '  -1.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0809'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D48").ClearFormats()

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.242.74'
$ws.Range("E49").Value = '  +8.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '2.21'
$ws.Range("E50").Value = '  -6.80%  '
$ws.Range("D50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.77'
$ws.Range("E51").Value = '  +12.67%  '
$ws.Range("D51").ClearFormats()

